$d = $word.ActiveDocument

# 1. Remove the leading "_GoBack" bookmark pair at the very start of the body.
#    The bookmark is zero-width and sits structurally before offset 0; an
#    InsertXML("") on the collapsed range consumes one of the two markers
#    (bookmarkStart, then bookmarkEnd) per call, so do it twice.
$bmRange = $d.Range(0, 0)
$bmRange.InsertXML("")
$bmRange2 = $d.Range(0, 0)
$bmRange2.InsertXML("")

# 2. First paragraph: bump the three ruby runs' hpsRaise 18 -> 46 and the
#    ruby text (w:rt) run size 48 -> 10, matching the other paragraphs.
$para1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000C40DD" w:rsidRPr="000C40DD" w:rsidRDefault="000C40DD" w:rsidP="000C40DD"><w:pPr><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r w:rsidRPr="000C40DD"><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:ruby><w:rubyPr><w:rubyAlign w:val="distributeSpace"/><w:hps w:val="10"/><w:hpsRaise w:val="46"/><w:hpsBaseText w:val="48"/><w:lid w:val="ja-JP"/></w:rubyPr><w:rt><w:r w:rsidR="000C40DD" w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:ascii="MS Mincho" w:hAnsi="MS Mincho" w:hint="eastAsia"/><w:sz w:val="10"/><w:szCs w:val="48"/></w:rPr><w:t>きもん</w:t></w:r></w:rt><w:rubyBase><w:r w:rsidR="000C40DD" w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>鬼門</w:t></w:r></w:rubyBase></w:ruby></w:r><w:r w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>の</w:t></w:r><w:r w:rsidRPr="000C40DD"><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:ruby><w:rubyPr><w:rubyAlign w:val="distributeSpace"/><w:hps w:val="10"/><w:hpsRaise w:val="46"/><w:hpsBaseText w:val="48"/><w:lid w:val="ja-JP"/></w:rubyPr><w:rt><w:r w:rsidR="000C40DD" w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:ascii="MS Mincho" w:hAnsi="MS Mincho" w:hint="eastAsia"/><w:sz w:val="10"/><w:szCs w:val="48"/></w:rPr><w:t>ほうがく</w:t></w:r></w:rt><w:rubyBase><w:r w:rsidR="000C40DD" w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>方角</w:t></w:r></w:rubyBase></w:ruby></w:r><w:r w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>を</w:t></w:r><w:r w:rsidRPr="000C40DD"><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:ruby><w:rubyPr><w:rubyAlign w:val="distributeSpace"/><w:hps w:val="10"/><w:hpsRaise w:val="46"/><w:hpsBaseText w:val="48"/><w:lid w:val="ja-JP"/></w:rubyPr><w:rt><w:r w:rsidR="000C40DD" w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:ascii="MS Mincho" w:hAnsi="MS Mincho" w:hint="eastAsia"/><w:sz w:val="10"/><w:szCs w:val="48"/></w:rPr><w:t>ぎょうし</w:t></w:r></w:rt><w:rubyBase><w:r w:rsidR="000C40DD" w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>凝視</w:t></w:r></w:rubyBase></w:ruby></w:r><w:r w:rsidRPr="000C40DD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>する。</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(1).Range.InsertXML($para1Xml)

# 3. Last (empty) paragraph: add a rightVertical ruby run ("中" / "ㄓㄨㄥ")
#    and move the "_GoBack" bookmark so it now sits after that new run.
$para7Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000C40DD" w:rsidRPr="00235942" w:rsidRDefault="000C40DD" w:rsidP="000C40DD"><w:pPr><w:rPr><w:rFonts w:eastAsia="新細明體"/><w:lang w:eastAsia="zh-TW"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="新細明體"/><w:lang w:eastAsia="zh-TW"/></w:rPr><w:ruby><w:rubyPr><w:rubyAlign w:val="rightVertical"/><w:hps w:val="7"/><w:hpsRaise w:val="18"/><w:hpsBaseText w:val="21"/><w:lid w:val="zh-TW"/></w:rubyPr><w:rt><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:hint="eastAsia"/><w:w w:val="75"/><w:sz w:val="7"/><w:lang w:eastAsia="zh-TW"/></w:rPr><w:t>ㄓㄨㄥ</w:t></w:r></w:rt><w:rubyBase><w:r><w:rPr><w:rFonts w:eastAsia="新細明體" w:hint="eastAsia"/><w:lang w:eastAsia="zh-TW"/></w:rPr><w:t>中</w:t></w:r></w:rubyBase></w:ruby></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$lastParaIndex = $d.Paragraphs.Count
$d.Paragraphs.Item($lastParaIndex).Range.InsertXML($para7Xml)
